# Fix shark double counts: update percent coverage values that were
# recalculated after correcting shark landings double counting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Global)
$ws.Range("C4").Value = 57.49218737885558
$ws.Range("D4").Value = 19.21686350833784
$ws.Range("E4").Value = 9.751189517155115
$ws.Range("G4").Value = 88.65776348812088

# Row 11 (Area 47)
$ws.Range("C11").Value = 74.29039347307072
$ws.Range("G11").Value = 95.38283744105355

# Row 12 (Area 51)
$ws.Range("D12").Value = 25.89026358993192
$ws.Range("G12").Value = 95.35980609945506

# Row 13 (Area 57)
$ws.Range("D13").Value = 34.8998241561326
$ws.Range("G13").Value = 98.48716029324659

# Row 16 (Area 71)
$ws.Range("E16").Value = 43.52029544476464
$ws.Range("G16").Value = 97.99543794616298

# Row 17 (Area 77)
$ws.Range("D17").Value = 9.888871821812476
$ws.Range("G17").Value = 87.65426515699602
